$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Edeq (row 4) phone numbers become text values
$ws.Range("B4").Value = "7558829"
$ws.Range("C4").Value = "75398"
$ws.Range("D4").Value = "3698547"

# Movistar (row 7) phone numbers become text values
$ws.Range("B7").Value = "3108228425"
$ws.Range("C7").Value = "3121715639"
$ws.Range("D7").Value = "312321666"

# Update the active selection from F12 to A6
$ws.Range("A6").Select()
